# billing/Samples/rates.xlsx -- "Rates excel and database"
#
# Rename the header row to the snake_case / lower-case column names used by
# the database import, and leave the cursor on D1 (next empty column) the
# way the workbook was left when it was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "product_id"
$ws.Range("B1").Value = "rate"
$ws.Range("C1").Value = "scope"

# Refresh the header/footer font label (some re-saves normalise
# "...,Regular" to "...,Normal") while keeping the same font/size.
$ps = $ws.PageSetup
$ps.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ps.CenterFooter = '&"Times New Roman,Normal"&12Page &P'

# Move the active selection to D1.
$ws.Range("D1").Select() | Out-Null
